$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E-column values (column 5) for the affected rows
$ws.Range("E2").Value  = 245150
$ws.Range("E3").Value  = 267755
$ws.Range("E4").Value  = 1510020
$ws.Range("E5").Value  = 2285967
$ws.Range("E6").Value  = 783876
$ws.Range("E7").Value  = 97630
$ws.Range("E8").Value  = 202700
$ws.Range("E9").Value  = 29320
$ws.Range("E10").Value = 24760
$ws.Range("E11").Value = 71508
$ws.Range("E12").Value = 56050
$ws.Range("E13").Value = 5277600
$ws.Range("E14").Value = 424800
$ws.Range("E19").Value = 386400
$ws.Range("E20").Value = 315600

# Update the active selection from E21 to E12
$ws.Range("E12").Select()
